# Update crypto price/volume data as scraped on Tue Jan  2 03:10:07 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.378.40"
$ws.Range("E2").Value = "  +5.75%  "

$ws.Range("D3").Value = "2.376.24"
$ws.Range("E3").Value = "  +3.56%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.41"
$ws.Range("E5").Value = "  +6.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "317.60"
$ws.Range("E6").Value = "  +1.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  +2.22%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +4.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.04"
$ws.Range("E10").Value = "  +7.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  +2.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.67"
$ws.Range("E12").Value = "  +4.74%  "

$ws.Range("E13").Value = "  +3.18%  "

$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.80"
$ws.Range("E15").Value = "  +4.61%  "

$ws.Range("D16").Value = "2.740.17"
$ws.Range("E16").Value = "  +3.69%  "

$ws.Range("D17").Value = "2.389.86"
$ws.Range("E17").Value = "  +3.90%  "

$ws.Range("D18").Value = "45.222.33"
$ws.Range("E18").Value = "  +6.34%  "

$ws.Range("E19").Value = "  +4.43%  "

$ws.Range("E20").Value = "  +3.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -3.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.11"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.54"
$ws.Range("E23").Value = "  +2.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.71"
$ws.Range("E24").Value = "  +1.71%  "

$ws.Range("E25").Value = "  +7.00%  "

$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.64"
$ws.Range("E27").Value = "  +8.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.28"
$ws.Range("E28").Value = "  +5.25%  "

$ws.Range("E29").Value = "  +2.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.92"
$ws.Range("E30").Value = "  +9.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.88"
$ws.Range("E31").Value = "  +2.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0934"
$ws.Range("E32").Value = "  +8.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "169.60"
$ws.Range("E33").Value = "  +2.83%  "

$ws.Range("E34").Value = "  +15.70%  "

$ws.Range("E35").Value = "  +1.93%  "

$ws.Range("E36").Value = "  +4.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.81"
$ws.Range("E37").Value = "  +6.85%  "

$ws.Range("E38").Value = "  +12.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0366"
$ws.Range("E39").Value = "  +4.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.90"
$ws.Range("E40").Value = "  +4.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  +9.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.03"
$ws.Range("E42").Value = "  +5.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.85"
$ws.Range("E43").Value = "  +15.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.240"
$ws.Range("E44").Value = "  +6.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.69"
$ws.Range("E45").Value = "  +3.63%  "

$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.73"
$ws.Range("E47").Value = "  +7.24%  "

$ws.Range("E48").Value = "  +10.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "79.62"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.61"
$ws.Range("E50").Value = "  +17.33%  "

# Row 51: coin changed from TheGraph to FraxShare
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.10"
$ws.Range("E51").Value = "  +5.13%  "
